$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 474.83334
$ws.Range("J33").Value = 400
$ws.Range("L33").Value = 400
$ws.Range("N33").Value = -858
$ws.Range("H34").Value = 15848
$ws.Range("I34").Value = 15848
$ws.Range("K34").Value = 15848
$ws.Range("M34").Value = -15645
$ws.Range("H36").Value = 15848
$ws.Range("I36").Value = 15848
$ws.Range("K36").Value = 15848
$ws.Range("M36").Value = -15133
$ws.Range("H75").Value = 53657
$ws.Range("J75").Value = 53657
$ws.Range("L75").Value = 53657
$ws.Range("N75").Value = -55529
$ws.Range("H78").Value = 53657
$ws.Range("J78").Value = 53657
$ws.Range("L78").Value = 160971
$ws.Range("N78").Value = -170331
$ws.Range("H107").Value = 36638.38
$ws.Range("I107").Value = 45862.477
$ws.Range("J107").Value = 1279.3334
$ws.Range("K107").Value = 45862.477
$ws.Range("L107").Value = 1279.3334
$ws.Range("M107").Value = -43942.477
$ws.Range("N107").Value = -5119.3334
$ws.Range("H132").Value = 4120.057
$ws.Range("I132").Value = 4051.5757
$ws.Range("K132").Value = 12154.7271
$ws.Range("M132").Value = -9624.7271
$ws.Range("H138").Value = 6699.0347
$ws.Range("J138").Value = 8012.5
$ws.Range("L138").Value = 24037.5
$ws.Range("N138").Value = -34317.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5778.4
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 5778.4
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 5778.4
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -6004.4
$ws.Range("H5").Value = 548.8182
$ws.Range("I5").Value = 80.25
$ws.Range("K5").Value = 80.25
$ws.Range("M5").Value = 31.75
$ws.Range("H61").Value = 4456.36
$ws.Range("I61").Value = 2661.6667
$ws.Range("K61").Value = 2661.6667
$ws.Range("M61").Value = -2449.6667
$ws.Range("H110").Value = 149440.83
$ws.Range("I110").Value = 163363.73
$ws.Range("K110").Value = 163363.73
$ws.Range("M110").Value = -161318.73
$ws.Range("H116").Value = 5778.4
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 5778.4
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 5778.4
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -10366.4
$ws.Range("H136").Value = 4456.36
$ws.Range("I136").Value = 2661.6667
$ws.Range("K136").Value = 7985.000100000001
$ws.Range("M136").Value = -5435.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5778.4
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 5778.4
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 5778.4
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -6006.4
$ws.Range("H4").Value = 548.8182
$ws.Range("I4").Value = 80.25
$ws.Range("K4").Value = 80.25
$ws.Range("M4").Value = 34.75
$ws.Range("H105").Value = 3147.12
$ws.Range("I105").Value = 2293.5
$ws.Range("K105").Value = 2293.5
$ws.Range("M105").Value = -546.5
$ws.Range("H134").Value = 8938
$ws.Range("I134").Value = 8907
$ws.Range("K134").Value = 26721
$ws.Range("M134").Value = -24186

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4004.3333
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 4004.3333
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 4004.3333
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -4578.3333
$ws.Range("H99").Value = 6074.5
$ws.Range("I99").Value = 4636.4
$ws.Range("K99").Value = 4636.4
$ws.Range("M99").Value = -3138.4
$ws.Range("H113").Value = 4004.3333
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 4004.3333
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 4004.3333
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -8344.3333
$ws.Range("H126").Value = 6074.5
$ws.Range("I126").Value = 4636.4
$ws.Range("K126").Value = 13909.2
$ws.Range("M126").Value = -11439.2
$ws.Range("H134").Value = 4584.8887
$ws.Range("I134").Value = 4031.5833
$ws.Range("K134").Value = 12094.7499
$ws.Range("M134").Value = -9559.749899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 2500521.2
$ws.Range("I8").Value = 2500521.2
$ws.Range("K8").Value = 7501563.600000001
$ws.Range("M8").Value = -7501424.600000001
$ws.Range("H14").Value = 6825
$ws.Range("I14").Value = 6825
$ws.Range("K14").Value = 20475
$ws.Range("M14").Value = -20302

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 2008.1666
$ws.Range("I99").Value = 1409.8
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 1409.8
$ws.Range("L99").Value = 5000
$ws.Range("M99").Value = 836.2
$ws.Range("N99").Value = -9492
$ws.Range("H113").Value = 391218.97
$ws.Range("I113").Value = 563293.9399999999
$ws.Range("J113").Value = 4050.25
$ws.Range("K113").Value = 563293.9399999999
$ws.Range("L113").Value = 4050.25
$ws.Range("M113").Value = -561123.9399999999
$ws.Range("N113").Value = -8390.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1008597.5
$ws.Range("I7").Value = 2509244.2
$ws.Range("J7").Value = 8166.3335
$ws.Range("K7").Value = 2509244.2
$ws.Range("L7").Value = 8166.3335
$ws.Range("M7").Value = -2509132.2
$ws.Range("N7").Value = -8390.333500000001
$ws.Range("H46").Value = 5020.778
$ws.Range("I46").Value = 3887.4
$ws.Range("J46").Value = 6437.5
$ws.Range("K46").Value = 3887.4
$ws.Range("L46").Value = 6437.5
$ws.Range("M46").Value = -3699.4
$ws.Range("N46").Value = -6813.5
$ws.Range("H61").Value = 7586.6
$ws.Range("I61").Value = 5465.5713
$ws.Range("K61").Value = 5465.5713
$ws.Range("M61").Value = -5263.5713
$ws.Range("H82").Value = 5584.3335
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("H85").Value = 5584.3335
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("H113").Value = 7586.6
$ws.Range("I113").Value = 5465.5713
$ws.Range("K113").Value = 5465.5713
$ws.Range("M113").Value = -3295.5713
$ws.Range("H126").Value = 1008597.5
$ws.Range("I126").Value = 2509244.2
$ws.Range("J126").Value = 8166.3335
$ws.Range("K126").Value = 7527732.600000001
$ws.Range("L126").Value = 24499.0005
$ws.Range("M126").Value = -7525262.600000001
$ws.Range("N126").Value = -29439.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 17168.066
$ws.Range("I81").Value = 17168.066
$ws.Range("K81").Value = 34336.132
$ws.Range("M81").Value = -33275.132
$ws.Range("H84").Value = 17168.066
$ws.Range("I84").Value = 17168.066
$ws.Range("K84").Value = 171680.66
$ws.Range("M84").Value = -166376.66
$ws.Range("H113").Value = 1144.1154
$ws.Range("I113").Value = 1042.5555
$ws.Range("K113").Value = 3127.6665
$ws.Range("M113").Value = -957.6664999999998
$ws.Range("H132").Value = 3679.1562
$ws.Range("I132").Value = 2596.875
$ws.Range("K132").Value = 7790.625
$ws.Range("M132").Value = -5260.625

Write-Output "Applied all Jenova_Profits cell updates"